$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $ws.Range($cellRange).Formula = '="' + $escaped + '"'
    $ws.Range($cellRange).Copy($ws.Range($cellRange)) | Out-Null
    $ws.Range($cellRange).PasteSpecial(-4163) | Out-Null
}

Set-TextValue "D2" '67.019.23'
Set-TextValue "E2" '  +1.27%  '

Set-TextValue "D3" '3.116.32'
Set-TextValue "E3" '  +1.63%  '

Set-TextValue "E4" '  -0.02%  '

Set-TextValue "D5" '575.65'
Set-TextValue "E5" '  -0.45%  '

Set-TextValue "D6" '174.08'
Set-TextValue "E6" '  +3.83%  '

Set-TextValue "D7" '1.00'
Set-TextValue "E7" '  -0.07%  '

Set-TextValue "D8" '3.113.79'
Set-TextValue "E8" '  +1.56%  '

Set-TextValue "D9" '0.521'
Set-TextValue "E9" '  -0.18%  '

Set-TextValue "E10" '  -3.46%  '

Set-TextValue "D11" '0.153'
Set-TextValue "E11" '  +0.55%  '

Set-TextValue "D12" '0.479'
Set-TextValue "E12" '  -1.05%  '

Set-TextValue "D13" '0.0000247'
Set-TextValue "E13" '  -0.68%  '

Set-TextValue "D14" '37.18'
Set-TextValue "E14" '  +1.43%  '

Set-TextValue "E15" '  -0.98%  '

Set-TextValue "D16" '3.635.67'
Set-TextValue "E16" '  +1.47%  '

Set-TextValue "D17" '67.050.33'
Set-TextValue "E17" '  +1.10%  '

Set-TextValue "D18" '7.11'
Set-TextValue "E18" '  -0.76%  '

Set-TextValue "D19" '3.122.03'
Set-TextValue "E19" '  +1.64%  '

Set-TextValue "D20" '16.20'
Set-TextValue "E20" '  -0.20%  '

Set-TextValue "D21" '476.70'
Set-TextValue "E21" '  +3.02%  '

Set-TextValue "D22" '0.711'
Set-TextValue "E22" '  +0.10%  '

Set-TextValue "E23" '  +3.92%  '

Set-TextValue "B24" 'Litecoin'
Set-TextValue "C24" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D24" '83.84'
Set-TextValue "E24" '  +0.88%  '

Set-TextValue "B25" 'InternetComputer(DFINITY)'
Set-TextValue "C25" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D25" '13.29'
Set-TextValue "E25" '  +3.78%  '

Set-TextValue "D26" '2.30'
Set-TextValue "E26" '  +1.29%  '

Set-TextValue "B27" 'Dai'
Set-TextValue "C27" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D27" '1.00'
Set-TextValue "E27" '  +0.02%  '

Set-TextValue "B28" 'RenderToken'
Set-TextValue "C28" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D28" '9.99'
Set-TextValue "E28" '  -0.22%  '

Set-TextValue "D29" '2.42'
Set-TextValue "E29" '  +0.38%  '

Set-TextValue "D30" '7.96'
Set-TextValue "E30" '  -1.43%  '

Set-TextValue "D31" '2.66'
Set-TextValue "E31" '  -0.18%  '

Set-TextValue "D32" '28.61'
Set-TextValue "E32" '  +1.42%  '

Set-TextValue "D33" '0.0₃0968'
Set-TextValue "E33" '  -4.53%  '

Set-TextValue "E34" '  -2.42%  '

Set-TextValue "D35" '1.00'
Set-TextValue "E35" '  -0.18%  '

Set-TextValue "D36" '5.86'
Set-TextValue "E36" '  +0.02%  '

Set-TextValue "E37" '  -1.40%  '

Set-TextValue "D38" '47.64'
Set-TextValue "E38" '  -1.35%  '

Set-TextValue "E39" '  +2.61%  '

Set-TextValue "D40" '50.07'
Set-TextValue "E40" '  +0.18%  '

Set-TextValue "D41" '0.310'
Set-TextValue "E41" '  -0.92%  '

Set-TextValue "E42" '  +0.67%  '

Set-TextValue "D43" '8.60'
Set-TextValue "E43" '  -0.19%  '

Set-TextValue "D44" '2.807.16'
Set-TextValue "E44" '  +1.69%  '

Set-TextValue "D45" '0.0356'
Set-TextValue "E45" '  -1.24%  '

Set-TextValue "D46" '380.01'
Set-TextValue "E46" '  -0.23%  '

Set-TextValue "D47" '2.56'
Set-TextValue "E47" '  -10.84%  '

Set-TextValue "D48" '136.12'
Set-TextValue "E48" '  +1.84%  '

Set-TextValue "D50" '24.76'
Set-TextValue "E50" '  +1.38%  '

Set-TextValue "D51" '2.20'
Set-TextValue "E51" '  -0.49%  '
